# Apply the edit described by the diff: fill in the missing POM_ITERATION
# values in column C for rows 2, 3 and 8 on the "Login" sheet, and update
# the sheet's view/selection (scroll back to the top, select A2 instead of
# the previously-selected A16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Fill in the new POM_ITERATION (column C) values that were missing.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C8").Value = 1

# Reset the scroll position (drop topLeftCell="A5") and move the
# selection to A2 (was A16).
$ws.Range("A2").Select()
